$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,1,16.782487,50.347461,0.9252099721531751,0.9252099721531751,1,0.3333333333333333,0.01569233333333333,0.047077,0.03693539111407157,0.03693539111407157,0.2633563801663333,2.370207421497,0.03417299218411679,0.03417299218411679),
    @(3,1,16.782487,50.347461,0.9252099721531751,0.9252099721531751,3,1,0.3230143333333333,0.9690430000000001,0.7602859615386125,0.7602859615386125,5.420983849980334,48.78885464982301,0.7034241533035896,0.7034241533035896),
    @(3,1,16.782487,50.347461,0.9252099721531751,0.9252099721531751,2,0.6666666666666666,0.08615233333333333,0.258457,0.202778647347316,0.202778647347316,1.445850414186333,13.012653727677,0.1876128266654688,0.1876128266654687),
    @(3,1,1.091026,3.273078,0.0601477084462148,0.0601477084462148,1,0.3333333333333333,0.01569233333333333,0.047077,0.03693539111407157,0.03693539111407157,0.01712074366733333,0.154086693006,0.00222157913607609,0.00222157913607609),
    @(3,1,1.091026,3.273078,0.0601477084462148,0.0601477084462148,3,1,0.3230143333333333,0.9690430000000001,0.7602859615386125,0.7602859615386125,0.3524170360393334,3.171753324354,0.04572945835037454,0.04572945835037454),
    @(3,1,1.091026,3.273078,0.0601477084462148,0.0601477084462148,2,0.6666666666666666,0.08615233333333333,0.258457,0.202778647347316,0.202778647347316,0.09399443562733334,0.8459499206459999,0.01219667095976417,0.01219667095976417),
    @(2,0.6666666666666666,0.2655986666666667,0.7967960000000001,0.01464231940061012,0.01464231940061012,1,0.3333333333333333,0.01569233333333333,0.047077,0.03693539111407157,0.03693539111407157,0.004167862810222222,0.037510765292,0.0005408197938786928,0.0005408197938786927),
    @(2,0.6666666666666666,0.2655986666666667,0.7967960000000001,0.01464231940061012,0.01464231940061012,3,1,0.3230143333333333,0.9690430000000001,0.7602859615386125,0.7602859615386125,0.08579217624755557,0.7721295862280001,0.01113234988464835,0.01113234988464834),
    @(2,0.6666666666666666,0.2655986666666667,0.7967960000000001,0.01464231940061012,0.01464231940061012,2,0.6666666666666666,0.08615233333333333,0.258457,0.202778647347316,0.202778647347316,0.02288194486355556,0.205937503772,0.002969149722083084,0.002969149722083082),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = 5 + $j
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
